# Localization-status report refresh ("Generate Report for Archive").
#
# The previously "Ready for handoff" rows have since moved on in the
# localization pipeline, so the report now shows them as "In Translation".
# This text shows up in:
#   - Overview sheet, columns E (zh-cn) & F (de-de), rows 2-3
#   - zh-cn sheet,   column C (Status), rows 2-3
#   - de-de sheet,   column C (Status), rows 2-3
#
# Because the new status text is shorter than the old one, the (previously
# uniformly widened) status columns are narrowed back down to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- zh-cn sheet: Status column (C) ---
foreach ($addr in @("C2", "C3")) {
    $cell = $zhcn.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- de-de sheet: Status column (C) ---
foreach ($addr in @("C2", "C3")) {
    $cell = $dede.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- Resize the status columns now that the text is shorter ---
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth   # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = $newWidth   # F: de-de
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # C: Status
$dede.Columns.Item(3).ColumnWidth = $newWidth        # C: Status
